# Update the "Buying Opportunity" / "support Zone" / "Short buildup" table
# on the active sheet: refresh the ticker lists in columns B/C (and the
# odd E2/F2/F3 cells), and extend the table from 13 data rows (A1:F13)
# down to 18 data rows (A1:F18), continuing the running index in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A keeps a zero-based running index, bold/centered/bordered to
# match the header row -- extend that formatting into the five brand new
# rows by copying the last existing data row's look-and-feel down before
# writing any new values into it.
$ws.Range("A13:F13").Copy($ws.Range("A14:F18"))

# New ticker data for columns B (Buying Opportunity) and C (support
# Zone), keyed by row number 2..18.
$data = @(
    @{ Row = 2;  B = "NSE:ASAL";        C = "NSE:AGROPHOS" },
    @{ Row = 3;  B = "NSE:DBSTOCKBRO";  C = "NSE:ANGELONE" },
    @{ Row = 4;  B = "NSE:DLINKINDIA";  C = "NSE:APCOTEXIND" },
    @{ Row = 5;  B = "NSE:DREAMFOLKS";  C = "NSE:DELTAMAGNT" },
    @{ Row = 6;  B = "NSE:ELECON";      C = "NSE:EMAMIREAL" },
    @{ Row = 7;  B = "NSE:EMKAY";       C = "NSE:GMDCLTD" },
    @{ Row = 8;  B = "NSE:GODREJIND";   C = "NSE:HERCULES" },
    @{ Row = 9;  B = "NSE:GULPOLY";     C = "NSE:HLVLTD" },
    @{ Row = 10; B = "NSE:GVPTECH";     C = "NSE:KECL" },
    @{ Row = 11; B = "NSE:HINDMOTORS";  C = "NSE:MIDHANI" },
    @{ Row = 12; B = "NSE:HONDAPOWER";  C = "NSE:NAZARA" },
    @{ Row = 13; B = "NSE:INDSWFTLAB";  C = "NSE:NDTV" },
    @{ Row = 14; B = "NSE:INDUSTOWER";  C = "NSE:PTL" },
    @{ Row = 15; B = "NSE:KAUSHALYA";   C = "NSE:RAMCOIND" },
    @{ Row = 16; B = "NSE:MAFANG";      C = "NSE:RTNINDIA" },
    @{ Row = 17; B = "NSE:MAHEPC";      C = $null },
    @{ Row = 18; B = "NSE:MODIRUBBER";  C = $null }
)

foreach ($entry in $data) {
    $r = $entry.Row

    # Column A: running 0-based index.
    $ws.Cells.Item($r, 1).Value = $r - 2

    # Column B: Buying Opportunity ticker.
    $ws.Cells.Item($r, 2).Value = $entry.B

    # Column C: support Zone ticker (blank for the two new tail rows).
    if ($entry.C) {
        $ws.Cells.Item($r, 3).Value = $entry.C
    }
}

# Column E: Short buildup -- row 2 gains a ticker that wasn't there before.
$ws.Cells.Item(2, 5).Value = "NSE:BATAINDIA"

# Column F: FII ENTERING -- the two tickers that used to live here are
# dropped in this revision, so clear them back out to blank cells.
$ws.Range("F2:F3").ClearContents()
